# Natmi following Dr Hou advice
# Update the Tfpi-Lrp1 LR-pair sheet with refreshed NATMI statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=4.288808333333333;  H=12.866425;  I=0.1993390083529519; J=0.1993390083529519; K=3; M=17.16653;          N=51.49959;   O=0.0560345397128279; P=0.0560345397128279; Q=73.62395691841665;  R=662.61561226575;    S=0.01116986957986922; T=0.01116986957986922 }
    3  = @{ E=3; G=4.288808333333333;  H=12.866425;  I=0.1993390083529519; J=0.1993390083529519; K=3; M=256.4443053333333;  N=769.332916; O=0.8370788162388805; P=0.8370788162388805; Q=1099.840473749477;  R=9898.564263745298;  S=0.1668624611423213;  T=0.1668624611423213  }
    4  = @{ E=3; G=4.288808333333333;  H=12.866425;  I=0.1993390083529519; J=0.1993390083529519; K=3; M=32.74538866666666;  N=98.236166;  O=0.1068866440482915; P=0.1068866440482915; Q=140.4386957918389;  R=1263.94826212655;   S=0.02130667763076138; T=0.02130667763076138 }
    5  = @{ E=3; G=11.10519266666667;  H=33.315578;  I=0.5161569185865866; J=0.5161569185865866; K=3; M=17.16653;          N=51.49959;   O=0.0560345397128279; P=0.0560345397128279; Q=190.6376230681133;  R=1715.73860761302;   S=0.02892261535259096; T=0.02892261535259096 }
    6  = @{ E=3; G=11.10519266666667;  H=33.315578;  I=0.5161569185865866; J=0.5161569185865866; K=3; M=256.4443053333333;  N=769.332916; O=0.8370788162388805; P=0.8370788162388805; Q=2847.863418996161;  R=25630.77077096545;  S=0.4320640224039681;  T=0.4320640224039681  }
    7  = @{ E=3; G=11.10519266666667;  H=33.315578;  I=0.5161569185865866; J=0.5161569185865866; K=3; M=32.74538866666666;  N=98.236166;  O=0.1068866440482915; P=0.1068866440482915; Q=363.6438500882164;  R=3272.794650793948;  S=0.05517028083002747; T=0.05517028083002747 }
    8  = @{ E=3; G=6.121147333333333;  H=18.363442;  I=0.2845040730604615; J=0.2845040730604615; K=3; M=17.16653;          N=51.49959;   O=0.0560345397128279; P=0.0560345397128279; Q=105.0788593320867;  R=945.7097339887799;  S=0.01594205478036772; T=0.01594205478036772 }
    9  = @{ E=3; G=6.121147333333333;  H=18.363442;  I=0.2845040730604615; J=0.2845040730604615; K=3; M=256.4443053333333;  N=769.332916; O=0.8370788162388805; P=0.8370788162388805; Q=1569.733375739652;  R=14127.60038165687;  S=0.2381523326925911;  T=0.2381523326925911  }
    10 = @{ E=3; G=6.121147333333333;  H=18.363442;  I=0.2845040730604615; J=0.2845040730604615; K=3; M=32.74538866666666;  N=98.236166;  O=0.1068866440482915; P=0.1068866440482915; Q=200.4393485159302;  R=1803.954136643372;  S=0.03040968558750267; T=0.03040968558750267 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
